# Generate Report for Handoff
# Update the "Latest Handoff Datetime" for the 37cbfd96-... row on the
# per-locale handoff status sheets, reflecting a fresh handoff pass.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D4").Value = "2016-03-08 04:52:20"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D4").Value = "2016-03-08 04:52:30"
